$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: measurement datetime changes on every data row
$ws.Range("A2:A11").Value = "2020-06-02 01:00"

# Helper: write a value into a cell while preserving its original text
# (inlineStr) cell type instead of letting Excel auto-coerce numeric-looking
# strings into numbers.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "E2" "0.025"
Set-TextValue "F2" "0.024"
Set-TextValue "G2" "18"
Set-TextValue "H2" "42"

# Row 3
Set-TextValue "C3" "0.003"
Set-TextValue "D3" "0.8"
Set-TextValue "E3" "0.026"
Set-TextValue "F3" "0.030"
Set-TextValue "G3" "48"
Set-TextValue "H3" "64"

# Row 4
Set-TextValue "C4" "0.003"
Set-TextValue "D4" "0.5"
Set-TextValue "E4" "0.023"
Set-TextValue "F4" "0.023"
Set-TextValue "G4" "24"
Set-TextValue "H4" "40"

# Row 5
Set-TextValue "D5" "0.5"
Set-TextValue "E5" "0.026"
Set-TextValue "F5" "0.025"
Set-TextValue "G5" "25"
Set-TextValue "H5" "51"

# Row 6
Set-TextValue "E6" "0.028"
Set-TextValue "F6" "0.023"
Set-TextValue "G6" "25"
Set-TextValue "H6" "47"

# Row 7
Set-TextValue "E7" "0.026"
Set-TextValue "F7" "0.017"
Set-TextValue "G7" "21"
Set-TextValue "H7" "43"

# Row 8
Set-TextValue "D8" "0.4"
Set-TextValue "E8" "0.041"
Set-TextValue "F8" "0.018"
Set-TextValue "G8" "17"
Set-TextValue "H8" "59"

# Row 9
Set-TextValue "C9" "0.002"
Set-TextValue "D9" "0.4"
Set-TextValue "E9" "0.029"
Set-TextValue "F9" "0.021"
Set-TextValue "G9" "20"
Set-TextValue "H9" "48"

# Row 10
Set-TextValue "D10" "0.4"
Set-TextValue "E10" "0.025"
Set-TextValue "F10" "0.034"
Set-TextValue "G10" "25"
Set-TextValue "H10" "56"

# Row 11
Set-TextValue "D11" "0.5"
Set-TextValue "E11" "0.025"
Set-TextValue "F11" "0.027"
Set-TextValue "G11" "22"
Set-TextValue "H11" "45"
